# Staircasing the stimulus duration: update the timing columns for the
# three "driving + lexical" rows and widen column B; also refresh the
# current selection to match where editing left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (training_driving_lexical) ---
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.5
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 20

# --- Row 3 (lexical_wo_driving_roboto) ---
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = 40
$ws.Range("G3").Value = 40

# --- Row 4 (lexical_wo_driving_neuefrutigerworld) ---
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.5
$ws.Range("F4").Value = 40
$ws.Range("G4").Value = 40

# Widen column B - ColumnWidth is in character units, OOXML stores the
# raw pixel-quantised width, so 29.16667 (character width) lands exactly
# on a stored width of 30.
$ws.Columns.Item(2).ColumnWidth = 29.16667

# Scroll the view so column B is the leftmost visible column, then move
# the selection from F4 to E4, as left at the end of the editing session.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("E4").Select() | Out-Null
